$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Range("A2").Value = 0.3787193298339841
$ws.Range("B2").Value = 0.750096321105957
$ws.Range("C2").Value = -2.375997304916381
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
